# Apartado 6.2.2 - RobertoTorralba
# Adds a new row ("Contenido" / AdBlock Plus description / empty cell)
# at the bottom of the comparison table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Append a new row; it inherits the table's 3 columns/widths and the
# "jc=left" paragraph formatting used by the rest of the table.
$newRow = $t.Rows.Add()

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Cell 1: bold "Contenido" -------------------------------------------
$cell1 = $newRow.Cells.Item(1)
[void]$cell1.Range.InsertXML("<w:p $wNs><w:pPr><w:jc w:val='left'/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Contenido</w:t></w:r></w:p>")

# --- Cell 2: "AdBlock Plus es mas completo ..." --------------------------
$cell2 = $newRow.Cells.Item(2)
$aTilde = [char]0x00E1
$cell2Xml = "<w:p $wNs><w:pPr><w:jc w:val='left'/></w:pPr>" + `
    "<w:proofErr w:type='spellStart'/><w:r><w:t>AdBlock</w:t></w:r><w:proofErr w:type='spellEnd'/>" + `
    "<w:r><w:t xml:space='preserve'> Plus es m${aTilde}s completo respecto al potencial de bloqueo de publicidad. </w:t></w:r>" + `
    "</w:p>"
[void]$cell2.Range.InsertXML($cell2Xml)
# InsertXML appends after the cell's existing (empty) default paragraph
# for every cell except the first one in a row, so drop that leftover
# empty paragraph, leaving only the text we just inserted.
[void]$cell2.Range.Paragraphs.Item(1).Range.Delete()

# --- Cell 3: empty, left justified --------------------------------------
$cell3 = $newRow.Cells.Item(3)
[void]$cell3.Range.InsertXML("<w:p $wNs><w:pPr><w:jc w:val='left'/></w:pPr></w:p>")
[void]$cell3.Range.Paragraphs.Item(1).Range.Delete()

Write-Output "Added 'Contenido' row to the table."
